$d = $word.ActiveDocument

# The "URL to GitHub Repository:" label is a single bold run at the top of
# the document. Position a collapsed range right after that text (still
# inside the same run/paragraph, before the paragraph mark) and type the
# space plus the repository URL there, so the new text inherits the
# existing bold character formatting.
$rng = $d.Content.Duplicate
$found = $rng.Find.Execute("URL to GitHub Repository:", $true, $false, $false, $false, $false,
                            $true, 1, $false, "", 0)

if ($found) {
    $rng.Collapse(0)
    $rng.InsertAfter(" https://github.com/MCSquaredTech/week08")
}
